$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet "Export" lists accounts sorted by descending Saldo. This edit:
#   1. Re-prices two accounts that moved from a (small) negative balance to
#      a small positive balance, which also moves their row to a new spot
#      further up the (descending) list:
#        - 004584517 CAIO : -63.49  -> 36.51   (moves up, between 37.28 and 35.77)
#        - 004332783 IRON : -90.94  ->  9.06   (moves up, between 9.8 and 8.68)
#   2. Removes the accounts that no longer belong on the sheet (most of the
#      negative-balance tail), keeping only the three largest negative
#      balances (JUNIO, DILSON, RODOLFO).
# ---------------------------------------------------------------------------

# --- Step 1: delete rows, bottom-to-top so earlier row numbers stay valid ---

# VINICIUS (-23864.81) through RICARDO (-73899.39)
$ws.Rows("284:294").Delete()

# RENAN (-13760.94), MATEUS (-14857.93)
$ws.Rows("280:281").Delete()

# old CAIO (-63.49) / old IRON (-90.94) through JULIO (-8823.32)
$ws.Rows("259:278").Delete()

# --- Step 2: insert the two re-priced rows at their new sorted position ---

# 004332783 IRON now 9.06 -> belongs right before 005077648 DUNAS (8.68)
$ws.Rows(214).Insert()
$ws.Cells.Item(214,1).NumberFormat = "@"
$ws.Cells.Item(214,1).Value = "004332783"
$ws.Cells.Item(214,2).Value = "IRON"
$ws.Cells.Item(214,3).Value = 9.06

# 004584517 CAIO now 36.51 -> belongs right before 004806286 VERA (35.77)
$ws.Rows(163).Insert()
$ws.Cells.Item(163,1).NumberFormat = "@"
$ws.Cells.Item(163,1).Value = "004584517"
$ws.Cells.Item(163,2).Value = "CAIO"
$ws.Cells.Item(163,3).Value = 36.51
